$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price values look like plain numbers need to be forced to
# text (NumberFormat "@") so Excel keeps them as strings like the source data,
# matching the workbook convention where every Price/Volume cell is inline text.
$textCells = @('D4', 'D5', 'D6', 'D7', 'D8', 'D9', 'D10', 'D11', 'D13', 'D14', 'D17', 'D19', 'D20', 'D21', 'D22', 'D23', 'D24', 'D25', 'D26', 'D28', 'D29', 'D30', 'D31', 'D32', 'D34', 'D35', 'D36', 'D38', 'D40', 'D41', 'D42', 'D43', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '29.039.01'
$ws.Range('E2').Value = '  -0.02%  '
$ws.Range('D3').Value = '1.834.93'
$ws.Range('E3').Value = '  +0.41%  '
$ws.Range('D4').Value = '0.9988'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '242.65'
$ws.Range('E5').Value = '  +0.54%  '
$ws.Range('D6').Value = '0.6277'
$ws.Range('E6').Value = '  -1.39%  '
$ws.Range('D7').Value = '0.9989'
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('D8').Value = '0.07587'
$ws.Range('E8').Value = '  +3.51%  '
$ws.Range('D9').Value = '0.2926'
$ws.Range('E9').Value = '  -0.31%  '
$ws.Range('D10').Value = '22.60'
$ws.Range('E10').Value = '  -0.80%  '
$ws.Range('D11').Value = '0.07740'
$ws.Range('E11').Value = '  +1.10%  '
$ws.Range('D12').Value = '1.836.05'
$ws.Range('E12').Value = '  +0.51%  '
$ws.Range('D13').Value = '4.966'
$ws.Range('E13').Value = '  -0.28%  '
$ws.Range('D14').Value = '0.6653'
$ws.Range('E14').Value = '  +0.33%  '
$ws.Range('E15').Value = '  +16.46%  '
$ws.Range('E16').Value = '  +1.63%  '
$ws.Range('D17').Value = '6.067'
$ws.Range('E17').Value = '  +0.24%  '
$ws.Range('D18').Value = '29.056.43'
$ws.Range('E18').Value = '  +0.54%  '
$ws.Range('D19').Value = '227.01'
$ws.Range('E19').Value = '  +1.43%  '
$ws.Range('D20').Value = '12.41'
$ws.Range('E20').Value = '  +0.18%  '
$ws.Range('D21').Value = '0.9995'
$ws.Range('E21').Value = '  -0.13%  '
$ws.Range('D22').Value = '7.218'
$ws.Range('E22').Value = '  +1.38%  '
$ws.Range('D23').Value = '0.9998'
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('D24').Value = '159.58'
$ws.Range('E24').Value = '  +1.10%  '
$ws.Range('D25').Value = '8.510'
$ws.Range('E25').Value = '  +0.56%  '
$ws.Range('D26').Value = '0.1385'
$ws.Range('E26').Value = '  +1.23%  '
$ws.Range('E27').Value = '  +0.44%  '
$ws.Range('D28').Value = '1.493'
$ws.Range('E28').Value = '  -0.77%  '
$ws.Range('D29').Value = '4.100'
$ws.Range('E29').Value = '  +0.25%  '
$ws.Range('D30').Value = '4.016'
$ws.Range('E30').Value = '  -0.15%  '
$ws.Range('D31').Value = '1.192'
$ws.Range('E31').Value = '  -0.93%  '
$ws.Range('D32').Value = '0.05252'
$ws.Range('E32').Value = '  -0.72%  '
$ws.Range('E33').Value = '  +0.53%  '
$ws.Range('D34').Value = '0.7362'
$ws.Range('E34').Value = '  -0.15%  '
$ws.Range('D35').Value = '1.138'
$ws.Range('E35').Value = '  -1.21%  '
$ws.Range('D36').Value = '2.677'
$ws.Range('E36').Value = '  +0.86%  '
$ws.Range('D37').Value = '1.242.62'
$ws.Range('E37').Value = '  -3.79%  '
$ws.Range('D38').Value = '2.761'
$ws.Range('E38').Value = '  +0.63%  '
$ws.Range('E39').Value = '  +0.12%  '
$ws.Range('D40').Value = '6.358'
$ws.Range('E40').Value = '  +1.10%  '
$ws.Range('D41').Value = '0.8981'
$ws.Range('E41').Value = '  +0.20%  '
$ws.Range('D42').Value = '0.9993'
$ws.Range('E42').Value = '  -0.06%  '
$ws.Range('D43').Value = '101.92'
$ws.Range('E43').Value = '  -0.61%  '
$ws.Range('D44').Value = '1.985.23'
$ws.Range('E44').Value = '  +0.46%  '
$ws.Range('D45').Value = '0.00000000124'
$ws.Range('E45').Value = '  +3.16%  '
$ws.Range('D46').Value = '64.24'
$ws.Range('E46').Value = '  +0.37%  '
$ws.Range('D47').Value = '0.5110'
$ws.Range('E47').Value = '  -0.50%  '
$ws.Range('D48').Value = '0.4040'
$ws.Range('E48').Value = '  +1.49%  '
$ws.Range('D49').Value = '8.878'
$ws.Range('E49').Value = '  +1.81%  '
$ws.Range('D50').Value = '0.05765'
$ws.Range('D51').Value = '6.704'
$ws.Range('E51').Value = '  +0.27%  '

# Restore default (Normal) styling on the forced cells so only the displayed
# text changes - no lingering custom number format on the cell itself.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
